$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends with a "Mean" column in AZ (right after "Run 49" in AY).
# A new "Run 50" column of results is being added before the "Mean" column, so
# insert a blank column at AZ; this shifts the existing "Mean" column to BA.
$ws.Range("AZ1:AZ14").EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("AZ1").Value = "Run 50"

# Per-row value for the new "Run 50" run (same result repeated on every row,
# consistent with the other run columns in this table) and the recomputed
# "Mean" value (now living in the shifted-over BA column) that accounts for it.
$newRunValue = 7278010408.414439
$newMeanValue = 7018616791.354816

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 52).Value = $newRunValue   # column AZ
    $ws.Cells.Item($row, 53).Value = $newMeanValue  # column BA
}
